$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 570, pushing the existing rows 570:629 down to 572:631.
$ws.Rows("570:571").Insert()

# New row 570 ("Primera" record dated 2022-07-27 / serial 44769)
$ws.Range("A570").Value = 9
$ws.Range("B570").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C570").Value = "Metropolitana"
$ws.Range("D570").Value = 44769
$ws.Range("E570").Value = 13
$ws.Range("F570").Value = 100112009
$ws.Range("G570").Value = "Acelga"
$ws.Range("H570").Value = "Sin especificar"
$ws.Range("I570").Value = "Primera"
$ws.Range("J570").Value = 52
$ws.Range("K570").Value = 16000
$ws.Range("L570").Value = 16000
$ws.Range("M570").Value = 16000
$ws.Range("N570").Value = '$/docena de atados'
$ws.Range("O570").Value = "Región Metropolitana"
$ws.Range("P570").Value = 5333
$ws.Range("Q570").Value = 3
$ws.Range("R570").Value = "Hortaliza"

# New row 571 ("Segunda" record dated 2022-07-27 / serial 44769)
$ws.Range("A571").Value = 9
$ws.Range("B571").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C571").Value = "Metropolitana"
$ws.Range("D571").Value = 44769
$ws.Range("E571").Value = 13
$ws.Range("F571").Value = 100112009
$ws.Range("G571").Value = "Acelga"
$ws.Range("H571").Value = "Sin especificar"
$ws.Range("I571").Value = "Segunda"
$ws.Range("J571").Value = 25
$ws.Range("K571").Value = 14000
$ws.Range("L571").Value = 14000
$ws.Range("M571").Value = 14000
$ws.Range("N571").Value = '$/docena de atados'
$ws.Range("O571").Value = "Región Metropolitana"
$ws.Range("P571").Value = 4667
$ws.Range("Q571").Value = 3
$ws.Range("R571").Value = "Hortaliza"

# Make sure the date cells carry the same number format as the other date cells in column D.
$ws.Range("D570:D571").NumberFormat = $ws.Range("D572").NumberFormat
